# feat: add 2022-Q1 data
#
# The existing "总计" (summary) sheet is repurposed into a new "2022-Q1"
# per-fund holdings sheet (it keeps its original sheetId), and a brand new
# "总计" sheet is appended after it, containing the old summary rows with a
# new 2022-Q1 row inserted at the top.

function Set-TextValue($range, [string]$text) {
    # Force a value to be stored as TEXT (even if it looks numeric, e.g. a
    # fund code with a leading zero, or a decimal that should stay a
    # string) instead of being auto-coerced to a number by the plain
    # Value setter.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Turn the old "总计" sheet into the new "2022-Q1" fund-holdings sheet
#    (keeps the original sheetId/relationship id).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Stamp the existing styled header/index cells onto the newly-needed
# columns/rows before we start overwriting values.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats
$q1.Range("A2").Copy()
$q1.Range("A3").PasteSpecial(-4122)      # xlPasteFormats

# Drop the old 2021-Q1 / 2020-Q4 rows (rows 4 & 5) - the new sheet only
# has 2 data rows.
$q1.Range("A4:H5").Clear()

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "010706"
$q1.Range("C2").Value = "景顺长城景骊成长混合型证券投资基金"
Set-TextValue $q1.Range("D2") "1.13"
Set-TextValue $q1.Range("E2") "93.50"
Set-TextValue $q1.Range("F2") "5.09"
Set-TextValue $q1.Range("G2") "0.0575"
$q1.Range("H2").Value = 10

# Row 3
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "260115"
$q1.Range("C3").Value = "景顺长城中小盘混合"
Set-TextValue $q1.Range("D3") "0.96"
Set-TextValue $q1.Range("E3") "94.00"
Set-TextValue $q1.Range("F3") "5.14"
Set-TextValue $q1.Range("G3") "0.0493"
$q1.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2) Build the new "总计" sheet right after "2022-Q1". Duplicate an
#    existing fully-formatted sheet (rather than Worksheets.Add(), which
#    creates a bare sheet) so sheetPr/pageMargins/etc. come along for
#    free, then wipe its content and rebuild it from scratch.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($null, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)
$total.Name = "总计"
$total.Cells.Clear()

$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats
$q1.Range("A2").Copy()
$total.Range("A2:A6").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# Row 2 - new 2022-Q1 entry
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11

# Row 3 - previously row 2 (2021-Q4)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 1.26

# Row 4 - previously row 3 (2021-Q3)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 1.57

# Row 5 - previously row 4 (2021-Q1)
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q1"
$total.Range("C5").Value = 9
$total.Range("D5").Value = 1.19

# Row 6 - previously row 5 (2020-Q4)
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2020-Q4"
$total.Range("C6").Value = 4
$total.Range("D6").Value = 0.42

# Restore the originally active sheet/tab (untouched by this edit).
$wb.Worksheets.Item("2020-Q4").Activate()
